# Regenerate the handback-status report: refresh the "Correspond Handoff
# Datetime" (column D) and "Correspond Handback DateTime" (column G)
# timestamps on row 3 of the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-13 15:50:26"
$wsZhCn.Range("G3").Value = "2016-01-13 15:51:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-13 15:50:36"
$wsDeDe.Range("G3").Value = "2016-01-13 15:51:29"
